$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)

# Pre-format column I (date strings) as Text so "2012-02-29" is not auto-converted to a date serial
$ws.Range("I1:I16").NumberFormat = "@"

# Header row (row 1)
$ws.Range("B1").Value = 'bank'
$ws.Range("C1").Value = 'deposit_type'
$ws.Range("D1").Value = 'currency'
$ws.Range("E1").Value = 'owner'
$ws.Range("F1").Value = 'total'
$ws.Range("G1").Value = 'property_category'
$ws.Range("H1").Value = 'category'
$ws.Range("I1").Value = 'date'
$ws.Range("J1").Value = 'legislator_name'
$ws.Range("K1").Value = 'legislator_id'
$ws.Range("L1").Value = 'source_file'
$ws.Range("M1").Value = 'index'

# Row 2
$ws.Range("A2").Value = 47
$ws.Range("B2").Value = '台新國際商業銀行建北分行'
$ws.Range("C2").Value = '綜合存款'
$ws.Range("D2").Value = '新臺幣'
$ws.Range("E2").Value = '蔡慧敏'
$ws.Range("F2").Value = 9603
$ws.Range("G2").Value = 'deposit'
$ws.Range("H2").Value = 'normal'
$ws.Range("I2").Value = '2012-02-29'
$ws.Range("J2").Value = '陳根德'
$ws.Range("K2").Value = 833
$ws.Range("L2").Value = 'tmp28cf1'
$ws.Range("M2").Value = 47

# Row 3
$ws.Range("A3").Value = 48
$ws.Range("B3").Value = '台新國際商業銀行建北分行'
$ws.Range("C3").Value = '活期存款'
$ws.Range("D3").Value = '美金'
$ws.Range("E3").Value = '蔡慧敏'
$ws.Range("F3").Value = 327
$ws.Range("G3").Value = 'deposit'
$ws.Range("H3").Value = 'normal'
$ws.Range("I3").Value = '2012-02-29'
$ws.Range("J3").Value = '陳根德'
$ws.Range("K3").Value = 833
$ws.Range("L3").Value = 'tmp28cf1'
$ws.Range("M3").Value = 48

# Row 4
$ws.Range("A4").Value = 49
$ws.Range("B4").Value = '花旗(台灣)商業銀行桃園'
$ws.Range("C4").Value = '活期存款'
$ws.Range("D4").Value = '新臺幣'
$ws.Range("E4").Value = '蔡慧敏'
$ws.Range("F4").Value = 16168
$ws.Range("G4").Value = 'deposit'
$ws.Range("H4").Value = 'normal'
$ws.Range("I4").Value = '2012-02-29'
$ws.Range("J4").Value = '陳根德'
$ws.Range("K4").Value = 833
$ws.Range("L4").Value = 'tmp28cf1'
$ws.Range("M4").Value = 49

# Row 5
$ws.Range("A5").Value = 50
$ws.Range("B5").Value = '花旗(台灣)商業銀行桃園'
$ws.Range("C5").Value = '活期存款'
$ws.Range("D5").Value = '美金'
$ws.Range("E5").Value = '蔡慧敏'
$ws.Range("F5").Value = 346035
$ws.Range("G5").Value = 'deposit'
$ws.Range("H5").Value = 'normal'
$ws.Range("I5").Value = '2012-02-29'
$ws.Range("J5").Value = '陳根德'
$ws.Range("K5").Value = 833
$ws.Range("L5").Value = 'tmp28cf1'
$ws.Range("M5").Value = 50

# Row 6
$ws.Range("A6").Value = 51
$ws.Range("B6").Value = '第一商業銀行北桃'
$ws.Range("C6").Value = '活期儲蓄存款'
$ws.Range("D6").Value = '新臺幣'
$ws.Range("E6").Value = '陳根德'
$ws.Range("F6").Value = 6785
$ws.Range("G6").Value = 'deposit'
$ws.Range("H6").Value = 'normal'
$ws.Range("I6").Value = '2012-02-29'
$ws.Range("J6").Value = '陳根德'
$ws.Range("K6").Value = 833
$ws.Range("L6").Value = 'tmp28cf1'
$ws.Range("M6").Value = 51

# Row 7
$ws.Range("A7").Value = 52
$ws.Range("B7").Value = '第一商業銀行北桃'
$ws.Range("C7").Value = '活期儲蓄存款'
$ws.Range("D7").Value = '新臺幣'
$ws.Range("E7").Value = '蔡慧敏'
$ws.Range("F7").Value = 231736
$ws.Range("G7").Value = 'deposit'
$ws.Range("H7").Value = 'normal'
$ws.Range("I7").Value = '2012-02-29'
$ws.Range("J7").Value = '陳根德'
$ws.Range("K7").Value = 833
$ws.Range("L7").Value = 'tmp28cf1'
$ws.Range("M7").Value = 52

# Row 8
$ws.Range("A8").Value = 53
$ws.Range("B8").Value = '第一商業銀行北桃'
$ws.Range("C8").Value = '支票存款'
$ws.Range("D8").Value = '新臺幣'
$ws.Range("E8").Value = '蔡慧敏'
$ws.Range("F8").Value = 56812
$ws.Range("G8").Value = 'deposit'
$ws.Range("H8").Value = 'normal'
$ws.Range("I8").Value = '2012-02-29'
$ws.Range("J8").Value = '陳根德'
$ws.Range("K8").Value = 833
$ws.Range("L8").Value = 'tmp28cf1'
$ws.Range("M8").Value = 53

# Row 9
$ws.Range("A9").Value = 54
$ws.Range("B9").Value = '合作金庫商業銀行桃圜'
$ws.Range("C9").Value = '活期儲蓄存款'
$ws.Range("D9").Value = '新臺幣'
$ws.Range("E9").Value = '陳根德'
$ws.Range("F9").Value = 653370
$ws.Range("G9").Value = 'deposit'
$ws.Range("H9").Value = 'normal'
$ws.Range("I9").Value = '2012-02-29'
$ws.Range("J9").Value = '陳根德'
$ws.Range("K9").Value = 833
$ws.Range("L9").Value = 'tmp28cf1'
$ws.Range("M9").Value = 54

# Row 10
$ws.Range("A10").Value = 55
$ws.Range("B10").Value = '合作金庫商業銀行桃園'
$ws.Range("C10").Value = '活期儲蓄存款'
$ws.Range("D10").Value = '新臺幣'
$ws.Range("E10").Value = '蔡慧敏'
$ws.Range("F10").Value = 1264
$ws.Range("G10").Value = 'deposit'
$ws.Range("H10").Value = 'normal'
$ws.Range("I10").Value = '2012-02-29'
$ws.Range("J10").Value = '陳根德'
$ws.Range("K10").Value = 833
$ws.Range("L10").Value = 'tmp28cf1'
$ws.Range("M10").Value = 55

# Row 11
$ws.Range("A11").Value = 56
$ws.Range("B11").Value = '合作金庫商業銀行桃圜'
$ws.Range("C11").Value = '活期存款'
$ws.Range("D11").Value = '美金'
$ws.Range("E11").Value = '蔡慧敏'
$ws.Range("F11").Value = 50
$ws.Range("G11").Value = 'deposit'
$ws.Range("H11").Value = 'normal'
$ws.Range("I11").Value = '2012-02-29'
$ws.Range("J11").Value = '陳根德'
$ws.Range("K11").Value = 833
$ws.Range("L11").Value = 'tmp28cf1'
$ws.Range("M11").Value = 56

# Row 12
$ws.Range("A12").Value = 57
$ws.Range("B12").Value = '日盛國際商業銀行北桃圜'
$ws.Range("C12").Value = '活期儲蓄存款'
$ws.Range("D12").Value = '新臺幣'
$ws.Range("E12").Value = '蔡慧敏'
$ws.Range("F12").Value = 935
$ws.Range("G12").Value = 'deposit'
$ws.Range("H12").Value = 'normal'
$ws.Range("I12").Value = '2012-02-29'
$ws.Range("J12").Value = '陳根德'
$ws.Range("K12").Value = 833
$ws.Range("L12").Value = 'tmp28cf1'
$ws.Range("M12").Value = 57

# Row 13
$ws.Range("A13").Value = 58
$ws.Range("B13").Value = 'H盛國際商業銀行北桃圜'
$ws.Range("C13").Value = '活期儲蓄存款'
$ws.Range("D13").Value = '新臺幣'
$ws.Range("E13").Value = '陳根德'
$ws.Range("F13").Value = 1822
$ws.Range("G13").Value = 'deposit'
$ws.Range("H13").Value = 'normal'
$ws.Range("I13").Value = '2012-02-29'
$ws.Range("J13").Value = '陳根德'
$ws.Range("K13").Value = 833
$ws.Range("L13").Value = 'tmp28cf1'
$ws.Range("M13").Value = 58

# Row 14
$ws.Range("A14").Value = 59
$ws.Range("B14").Value = '聯邦商業銀行桃圜'
$ws.Range("C14").Value = '綜合存款'
$ws.Range("D14").Value = '新臺幣'
$ws.Range("E14").Value = '蔡慧敏'
$ws.Range("F14").Value = 388230
$ws.Range("G14").Value = 'deposit'
$ws.Range("H14").Value = 'normal'
$ws.Range("I14").Value = '2012-02-29'
$ws.Range("J14").Value = '陳根德'
$ws.Range("K14").Value = 833
$ws.Range("L14").Value = 'tmp28cf1'
$ws.Range("M14").Value = 59

# Row 15
$ws.Range("A15").Value = 60
$ws.Range("B15").Value = '中華郵政股份有限公司桃園府前'
$ws.Range("C15").Value = '活期儲蓄存款'
$ws.Range("D15").Value = '新臺幣'
$ws.Range("E15").Value = '蔡慧敏'
$ws.Range("F15").Value = 3046
$ws.Range("G15").Value = 'deposit'
$ws.Range("H15").Value = 'normal'
$ws.Range("I15").Value = '2012-02-29'
$ws.Range("J15").Value = '陳根德'
$ws.Range("K15").Value = 833
$ws.Range("L15").Value = 'tmp28cf1'
$ws.Range("M15").Value = 60

# Row 16
$ws.Range("A16").Value = 61
$ws.Range("B16").Value = '中華郵政股份有限公司桃圜府前'
$ws.Range("C16").Value = '活期存款'
$ws.Range("D16").Value = '新臺幣'
$ws.Range("E16").Value = '蔡慧敏'
$ws.Range("F16").Value = 80517
$ws.Range("G16").Value = 'deposit'
$ws.Range("H16").Value = 'normal'
$ws.Range("I16").Value = '2012-02-29'
$ws.Range("J16").Value = '陳根德'
$ws.Range("K16").Value = 833
$ws.Range("L16").Value = 'tmp28cf1'
$ws.Range("M16").Value = 61

